$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.11903208268695039673
$ws.Range("C2").Value = 0.00006240767534437808
$ws.Range("D2").Value = 0.14942197473980470268
$ws.Range("E2").Value = 0.49423653606076972666
$ws.Range("G2").Value = 0.76275300116286914864

# Row 3
$ws.Range("B3").Value = 0.04271373187048221887
$ws.Range("C3").Value = 0.04071648406533733694
$ws.Range("D3").Value = 3.53776164880671917246
$ws.Range("E3").Value = 0.49423653606076972666
$ws.Range("G3").Value = 4.11542840080330840635
